{"js": "// FCREPO-1049: append an extra note about the writeLimit testing to the\n// test document, mirroring the commit that added this paragraph (plus the\n// blank paragraphs before/after it) right after the existing\n// \"This is a test document...\" paragraph and before the section break.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The target text always sits in the last paragraph of the body (just\n// before the final section properties), so anchor on it explicitly.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert a blank paragraph, then the new descriptive paragraph, then\n// another blank paragraph - all inserted after the last existing\n// paragraph, in document order.\nconst blankBefore = lastParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nconst newParagraph = blankBefore.insertParagraph(\n  \"Extension for testing fcrepo-1049 Improve the control over writeLimit in getDatastreamFromTika. The test sets writeLimit to a low number and reindexes this object, then search with gfindObjects will reveal that only the first writeLimit characters were used in indexing the datastream.\",\n  \"After\"\n);\nawait context.sync();\n\nnewParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n", "ps1": "# FCREPO-1049: append an extra note about the writeLimit testing to the\n# test document, mirroring the commit that added this paragraph (plus the\n# blank paragraphs before/after it) right after the existing\n# \"This is a test document...\" paragraph and before the section break.\n\n$d = $word.ActiveDocument\n\n# Step 1: insert a blank paragraph right after the current last paragraph\n# (\"This is a test document for the Tika extraction of Gsearch 2.4\n# fcrepo-1010.\").\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n# Step 2: insert a new (currently empty) paragraph after that blank one,\n# then fill it in with the new descriptive text.\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertBefore(\"Extension for testing fcrepo-1049 Improve the control over writeLimit in getDatastreamFromTika. The test sets writeLimit to a low number and reindexes this object, then search with gfindObjects will reveal that only the first writeLimit characters were used in indexing the datastream.\")\n\n# Step 3: insert a trailing blank paragraph after the new text paragraph.\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n$d.Save()\n"}
